$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stash a pristine copy of the hyperlink-cell formatting (currently style
# index 4 / "Hipervinculo") in an unused scratch cell so we can restore it
# later, after Hyperlinks.Add() re-applies its own (slightly different)
# hyperlink style to the target cell.
$ws.Range("B4").Copy()
$ws.Range("K1").PasteSpecial(-4122)

# Drop the two existing DataSource rows (old row 2 and row 4) but keep the
# formatting that lived on row 2 (row height, text styles) by removing rows
# 3 and 4 instead of row 2. This way the surviving row (old row 4) slides up
# into row 2 without losing any row/cell formatting.
$ws.Hyperlinks.Delete()
$ws.Range("A3:I4").EntireRow.Delete()

# Update row 2 with the new DataSource values
$ws.Range("A2").Value = "i-preproducciongestion.segurossura.com.ar"
$ws.Range("C2").Value = "su"
$ws.Range("D2").Value = "silverarrow"
$ws.Range("E2").Value = "'04104013020"
$ws.Range("B2").Value = "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"

# Re-create the hyperlink on B2 pointing at the new URL
$ws.Hyperlinks.Add($ws.Range("B2"), "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do")

# Restore the original hyperlink-cell formatting and drop the scratch cell
$ws.Range("K1").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("K1").Clear()

# Restore the selected cell shown when the sheet was last saved
$ws.Range("E3").Select()
